$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @('35.636.33', '  +0.41%  ')
    3 = @('1.897.95', '  +0.06%  ')
    4 = @('1.00', '  -0.03%  ')
    5 = @('247.35', '  -0.10%  ')
    6 = @('0.694', '  +0.01%  ')
    7 = @('1.00', '  +0.03%  ')
    8 = @('43.25', '  -1.60%  ')
    9 = @('57.58', '  +10.79%  ')
    10 = @('0.359', '  +2.10%  ')
    11 = @('0.0757', '  +2.08%  ')
    12 = @('0.0986', '  +1.56%  ')
    13 = @('14.58', '  +11.00%  ')
    14 = @('0.808', '  +11.59%  ')
    15 = @('2.173.06', '  +0.13%  ')
    16 = @('5.05', '  +2.63%  ')
    17 = @('1.922.04', '  +1.34%  ')
    18 = @('35.620.28', '  +0.56%  ')
    19 = @('73.90', '  +0.45%  ')
    20 = @('0.0₃0834', '  +1.29%  ')
    21 = @('247.18', '  -0.30%  ')
    22 = @('13.04', '  +1.53%  ')
    23 = @('5.22', '  +5.24%  ')
    24 = @('2.66', '  +5.35%  ')
    25 = @('1.00', '  +0.05%  ')
    26 = @('2.16', '  -1.12%  ')
    27 = @('166.74', '  +0.62%  ')
    28 = @('8.67', '  +2.11%  ')
    29 = @('18.40', '  +0.03%  ')
    30 = @('0.128', '  +0.75%  ')
    31 = @('4.39', '  +3.32%  ')
    32 = @('0.0602', '  +3.65%  ')
    33 = @('4.27', '  +0.46%  ')
    34 = @('1.84', '  +14.69%  ')
    35 = @('1.00', '  -0.05%  ')
    36 = @('1.48', '  -17.23%  ')
    37 = @('0.859', '  +0.59%  ')
    38 = @('1.97', '  -1.62%  ')
    39 = @('0.0731', '  +7.51%  ')
    40 = @('0.0227', '  +6.77%  ')
    41 = @('99.65', '  +2.33%  ')
    42 = @('17.17', '  -0.19%  ')
    43 = @('1.09', '  -0.24%  ')
    44 = @('14.29', '  +18.10%  ')
    45 = @('1.323.93', '  +2.89%  ')
    46 = @('2.38', '  +0.74%  ')
    47 = @('0.0811', '  +1.24%  ')
    48 = @('2.39', '  -0.48%  ')
    49 = @('2.74', '  -0.38%  ')
    50 = @('6.43', '  +0.61%  ')
    51 = @('42.80', '  -0.81%  ')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[0]
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $vals[1]
}